$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.422.42'
$ws.Range('E2').Value = '  -4.25%  '

$ws.Range('D3').Value = '2.539.94'
$ws.Range('E3').Value = '  -3.95%  '

$ws.Range('E4').Value = '  +0.07%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '505.81'
$ws.Range('E5').Value = '  -4.76%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '142.79'
$ws.Range('E6').Value = '  -8.18%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  +0.13%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.563'
$ws.Range('E8').Value = '  -4.80%  '

$ws.Range('D9').Value = '2.546.10'
$ws.Range('E9').Value = '  -4.14%  '

$ws.Range('E10').Value = '  -8.36%  '

$ws.Range('E11').Value = '  -6.91%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.331'
$ws.Range('E12').Value = '  -5.84%  '

$ws.Range('E13').Value = '  -0.54%  '

$ws.Range('D14').Value = '2.989.83'
$ws.Range('E14').Value = '  -3.75%  '

$ws.Range('D15').Value = '58.423.61'
$ws.Range('E15').Value = '  -4.25%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '20.67'
$ws.Range('E16').Value = '  -6.05%  '

$ws.Range('E17').Value = '  -6.24%  '

$ws.Range('D18').Value = '2.548.60'
$ws.Range('E18').Value = '  -6.67%  '

$ws.Range('E19').Value = '  -5.11%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '338.24'
$ws.Range('E20').Value = '  -4.55%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '10.06'
$ws.Range('E21').Value = '  -5.77%  '

$ws.Range('E22').Value = '  -0.05%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '5.94'
$ws.Range('E23').Value = '  -4.78%  '

$ws.Range('E24').Value = '  -1.87%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.410'
$ws.Range('E25').Value = '  -4.75%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.999'
$ws.Range('E26').Value = '  -0.31%  '

$ws.Range('D27').Value = '2.655.77'
$ws.Range('E27').Value = '  -3.77%  '

$ws.Range('E28').Value = '  -5.74%  '

$ws.Range('D29').Value = '0.0₃0785'
$ws.Range('E29').Value = '  -9.09%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '6.93'
$ws.Range('E30').Value = '  -6.20%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.00'
$ws.Range('E31').Value = '  +0.04%  '

$ws.Range('E32').Value = '  -0.74%  '

$ws.Range('B33').Value = 'EthereumClassic'
$ws.Range('C33').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '18.49'
$ws.Range('E33').Value = '  -5.23%  '

$ws.Range('B34').Value = 'Aptos'
$ws.Range('C34').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '5.82'
$ws.Range('E34').Value = '  -5.17%  '

$ws.Range('E35').Value = '  -5.99%  '

$ws.Range('B36').Value = 'NEARProtocol'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '3.88'
$ws.Range('E36').Value = '  -6.57%  '

$ws.Range('B37').Value = 'SuiNetwork'
$ws.Range('C37').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.906'
$ws.Range('E37').Value = '  +1.25%  '

$ws.Range('E38').Value = '  -8.01%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '35.94'
$ws.Range('E39').Value = '  -1.84%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.819'
$ws.Range('E40').Value = '  -11.42%  '

$ws.Range('E41').Value = '  -7.27%  '

$ws.Range('B42').Value = 'Filecoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '3.51'
$ws.Range('E42').Value = '  -8.10%  '

$ws.Range('B43').Value = 'Bittensor'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '280.00'
$ws.Range('E43').Value = '  -8.70%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.0995'
$ws.Range('E44').Value = '  -2.32%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.998'
$ws.Range('E45').Value = '  +0.08%  '

$ws.Range('E46').Value = '  -6.96%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0532'
$ws.Range('E47').Value = '  -5.44%  '

$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '18.71'
$ws.Range('E48').Value = '  -5.68%  '

$ws.Range('B49').Value = 'WhiteBITCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '10.29'
$ws.Range('E49').Value = '  -0.66%  '

$ws.Range('E50').Value = '  -5.42%  '

$ws.Range('E51').Value = '  -9.18%  '
